$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 403; existing rows 403-460 shift down to 404-461.
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row 403 with the new record.
$ws.Range("A403").Value = 11
$ws.Range("B403").Value = "Vega Monumental Concepción"
$ws.Range("C403").Value = "Bíobío"
$ws.Range("D403").Value = 45218
$ws.Range("E403").Value = 8
$ws.Range("F403").Value = 100112045
$ws.Range("G403").Value = "Zapallo"
$ws.Range("H403").Value = "Camote"
$ws.Range("I403").Value = "1a nueva(o)"
$ws.Range("J403").Value = 300
$ws.Range("K403").Value = 1200
$ws.Range("L403").Value = 1200
$ws.Range("M403").Value = 1200
$ws.Range("N403").Value = "$/kilo (volumen en unidades)"
$ws.Range("O403").Value = "Perú"
$ws.Range("P403").Value = 1200
$ws.Range("Q403").Value = 1
$ws.Range("R403").Value = "Hortaliza"
